# Insert a new weekly price record for Mandarina "Murcott" into the
# consolidated daily logic sheet. The new record is placed at row 32,
# pushing the existing rows 32-72 down to 33-73 (dimension grows to A1:T73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32 and below down by one row.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value() = 1
$ws.Cells.Item(32, 2).Value() = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value() = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value() = 44482
$ws.Cells.Item(32, 5).Value() = 15
$ws.Cells.Item(32, 6).Value() = "Fruta"
$ws.Cells.Item(32, 7).Value() = 100102
$ws.Cells.Item(32, 8).Value() = "Cítricos"
$ws.Cells.Item(32, 9).Value() = 100102004
$ws.Cells.Item(32, 10).Value() = "Mandarina"
$ws.Cells.Item(32, 11).Value() = "Murcott"
$ws.Cells.Item(32, 12).Value() = "Segunda"
$ws.Cells.Item(32, 13).Value() = 250
$ws.Cells.Item(32, 14).Value() = 12000
$ws.Cells.Item(32, 15).Value() = 13000
$ws.Cells.Item(32, 16).Value() = 12500
$ws.Cells.Item(32, 17).Value() = "$/caja 20 kilos"
$ws.Cells.Item(32, 18).Value() = "Región de Coquimbo"
$ws.Cells.Item(32, 19).Value() = 625
$ws.Cells.Item(32, 20).Value() = 20
